$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Femacal de La Calera - Coliflor"
# data block (before the former row 406), pushing the existing rows
# 406-429 down to 408-431. This reflects a new weekly price observation
# being added to the series.
$ws.Rows.Item(406).Insert()
$ws.Rows.Item(406).Insert()

# New row 406: Coliflor, Primera, fecha 44516 (2021-11-16)
$ws.Range("A406").Value = 3
$ws.Range("B406").Value = "Femacal de La Calera"
$ws.Range("C406").Value = "Coquimbo"
$ws.Range("D406").Value = 44516
$ws.Range("E406").Value = 5
$ws.Range("F406").Value = 100112008
$ws.Range("G406").Value = "Coliflor"
$ws.Range("H406").Value = "Sin especificar"
$ws.Range("I406").Value = "Primera"
$ws.Range("J406").Value = 850
$ws.Range("K406").Value = 600
$ws.Range("L406").Value = 600
$ws.Range("M406").Value = 600
$ws.Range("N406").Value = "$/unidad"
$ws.Range("O406").Value = "Provincia de Quillota"
$ws.Range("P406").Value = 600
$ws.Range("Q406").Value = 1
$ws.Range("R406").Value = "Hortaliza"

# New row 407: Coliflor, Segunda, fecha 44516 (2021-11-16)
$ws.Range("A407").Value = 3
$ws.Range("B407").Value = "Femacal de La Calera"
$ws.Range("C407").Value = "Coquimbo"
$ws.Range("D407").Value = 44516
$ws.Range("E407").Value = 5
$ws.Range("F407").Value = 100112008
$ws.Range("G407").Value = "Coliflor"
$ws.Range("H407").Value = "Sin especificar"
$ws.Range("I407").Value = "Segunda"
$ws.Range("J407").Value = 950
$ws.Range("K407").Value = 500
$ws.Range("L407").Value = 500
$ws.Range("M407").Value = 500
$ws.Range("N407").Value = "$/unidad"
$ws.Range("O407").Value = "Provincia de Quillota"
$ws.Range("P407").Value = 500
$ws.Range("Q407").Value = 1
$ws.Range("R407").Value = "Hortaliza"
